$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "35.371.97"
$ws.Range("E2").Value = "  -3.72%  "

# Row 3
$ws.Range("D3").Value = "1.979.43"
$ws.Range("E3").Value = "  -4.93%  "

# Row 4
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").Value = "'239.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "

# Row 6
$ws.Range("D6").Value = "'0.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -14.25%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").Value = "'55.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.22%  "

# Row 9
$ws.Range("D9").Value = "'59.01"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.355"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.16%  "

# Row 11
$ws.Range("D11").Value = "'0.0723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.84%  "

# Row 12
$ws.Range("D12").Value = "'0.102"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.55%  "

# Row 13
$ws.Range("D13").Value = "'0.892"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.20%  "

# Row 14
$ws.Range("D14").Value = "'14.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.97%  "

# Row 15
$ws.Range("D15").Value = "2.270.97"
$ws.Range("E15").Value = "  -5.16%  "

# Row 16
$ws.Range("E16").Value = "  -4.07%  "

# Row 17
$ws.Range("D17").Value = "1.978.97"
$ws.Range("E17").Value = "  -4.74%  "

# Row 18
$ws.Range("D18").Value = "'16.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.71%  "

# Row 19
$ws.Range("D19").Value = "35.224.35"
$ws.Range("E19").Value = "  -4.05%  "

# Row 20
$ws.Range("D20").Value = "'69.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.25%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0830"
$ws.Range("E21").Value = "  -4.83%  "

# Row 22
$ws.Range("D22").Value = "'230.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.30%  "

# Row 23
$ws.Range("D23").Value = "'4.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.69%  "

# Row 24
$ws.Range("E24").Value = "  +0.11%  "

# Row 25
$ws.Range("D25").Value = "'2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.07%  "

# Row 26
$ws.Range("D26").Value = "'2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.44%  "

# Row 27
$ws.Range("D27").Value = "'162.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.87%  "

# Row 28
$ws.Range("D28").Value = "'9.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.78%  "

# Row 29
$ws.Range("D29").Value = "'19.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.22%  "

# Row 30
$ws.Range("E30").Value = "  -12.57%  "

# Row 31
$ws.Range("D31").Value = "'1.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "

# Row 32
$ws.Range("D32").Value = "'4.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.44%  "

# Row 33
$ws.Range("D33").Value = "'0.0581"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.63%  "

# Row 34
$ws.Range("E34").Value = "  +8.16%  "

# Row 35
$ws.Range("E35").Value = "  -9.40%  "

# Row 36
$ws.Range("E36").Value = "  -0.26%  "

# Row 37
$ws.Range("D37").Value = "'2.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.27%  "

# Row 38
$ws.Range("E38").Value = "  -2.44%  "

# Row 39
$ws.Range("D39").Value = "'4.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
$ws.Range("E40").Value = "  -6.24%  "

# Row 41
$ws.Range("D41").Value = "'2.79"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "

# Row 42
$ws.Range("E42").Value = "  -5.76%  "

# Row 43
$ws.Range("E43").Value = "  -6.60%  "

# Row 44
$ws.Range("E44").Value = "  -8.88%  "

# Row 45
$ws.Range("D45").Value = "'89.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.86%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'7.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.01%  "

# Row 47
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.347.67"
$ws.Range("E47").Value = "  -3.22%  "

# Row 48
$ws.Range("D48").Value = "'15.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.29%  "

# Row 49
$ws.Range("E49").Value = "  -0.62%  "

# Row 50
$ws.Range("E50").Value = "  -6.99%  "

# Row 51
$ws.Range("D51").Value = "'45.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.49%  "
